# Auto-generated edit script applying the cryptos.xlsx data refresh described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds price text that LOOKS like a plain number (e.g. "1.00", "236.31").
# The sheet stores these as literal text (inline strings), not numeric values, so for those
# cells we pre-format as Text ("@") before assigning - otherwise Excel would silently convert
# "1.00" -> 1. Values with multiple dots (e.g. "96.489.96") are never number-like and are
# assigned directly.

# Row 2
$ws.Range("D2").Value = '96.489.96'
$ws.Range("E2").Value = '  -0.29%  '
# Row 3
$ws.Range("D3").Value = '3.695.11'
$ws.Range("E3").Value = '  +1.26%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.31'
$ws.Range("E5").Value = '  -2.48%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.88'
$ws.Range("E6").Value = '  +1.70%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '651.26'
$ws.Range("E7").Value = '  -0.71%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.430'
$ws.Range("E8").Value = '  +1.69%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '
# Row 10
$ws.Range("E10").Value = '  -0.88%  '
# Row 11
$ws.Range("D11").Value = '3.694.52'
$ws.Range("E11").Value = '  +1.35%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.32'
$ws.Range("E12").Value = '  +0.19%  '
# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000306'
$ws.Range("E13").Value = '  +18.24%  '
# Row 14
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.204'
$ws.Range("E14").Value = '  -0.03%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.74'
$ws.Range("E15").Value = '  +3.81%  '
# Row 16
$ws.Range("D16").Value = '4.382.61'
$ws.Range("E16").Value = '  +1.24%  '
# Row 17
$ws.Range("D17").Value = '96.257.97'
$ws.Range("E17").Value = '  -0.24%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.78'
$ws.Range("E18").Value = '  +13.04%  '
# Row 19
$ws.Range("D19").Value = '3.715.71'
$ws.Range("E19").Value = '  +1.74%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.92'
$ws.Range("E20").Value = '  +0.27%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.78'
$ws.Range("E21").Value = '  +2.66%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.503'
$ws.Range("E22").Value = '  -5.73%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '517.00'
$ws.Range("E23").Value = '  +0.94%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.37'
$ws.Range("E24").Value = '  -2.14%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000208'
$ws.Range("E25").Value = '  +1.47%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.94'
$ws.Range("E26").Value = '  +0.85%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '100.70'
$ws.Range("E27").Value = '  -0.42%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.06'
$ws.Range("E28").Value = '  +0.29%  '
# Row 29
$ws.Range("E29").Value = '  +1.82%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.00'
$ws.Range("E30").Value = '  -1.11%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.07'
$ws.Range("E31").Value = '  +1.99%  '
# Row 32
$ws.Range("E32").Value = '  -0.08%  '
# Row 33
$ws.Range("E33").Value = '  +5.43%  '
# Row 34
$ws.Range("E34").Value = '  -0.92%  '
# Row 35
$ws.Range("E35").Value = '  +0.42%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '653.31'
$ws.Range("E36").Value = '  +5.95%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '32.10'
$ws.Range("E37").Value = '  -3.02%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.585'
$ws.Range("E38").Value = '  +0.31%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.82'
$ws.Range("E39").Value = '  -0.08%  '
# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.86'
$ws.Range("E41").Value = '  +11.80%  '
# Row 42
$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.06'
$ws.Range("E42").Value = '  +6.45%  '
# Row 43
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.78'
$ws.Range("E43").Value = '  -3.45%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.159'
$ws.Range("E44").Value = '  +0.51%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.956'
$ws.Range("E45").Value = '  +0.44%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0454'
$ws.Range("E46").Value = '  +2.44%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.431'
$ws.Range("E47").Value = '  +4.04%  '
# Row 48
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.57'
$ws.Range("E48").Value = '  -0.17%  '
# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.27'
$ws.Range("E49").Value = '  -1.87%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.42'
$ws.Range("E50").Value = '  -2.03%  '
# Row 51
$ws.Range("E51").Value = '  +2.26%  '
